# Cascade workbook: add a "Cascade Characteristics" worksheet between
# "Compartments" and "Transitions", populate it with the two example
# characteristics (Latent / Active infections), and restore selections.

$wb = $excel.ActiveWorkbook

# --- Insert the new worksheet right before "Transitions" -------------------
$transitions = $wb.Worksheets.Item("Transitions")
$ws = $wb.Worksheets.Add($transitions)
$ws.Name = "Cascade Characteristics"

# --- Header row --------------------------------------------------------
$ws.Range("A1").Value = "Code Label"
$ws.Range("B1").Value = "Full Name"
$ws.Range("C1").Value = "Includes"

# --- Code labels (column A) ---------------------------------------------
$ws.Range("A2").Value = "lt_inf"
$ws.Range("A3").Value = "ac_inf"

# --- Full names (column B) ----------------------------------------------
$ws.Range("B2").Value = "Latent Infections"
$ws.Range("B3").Value = "Active Infections"

# --- Included compartments (columns C/D) --------------------------------
$ws.Range("C2").Value = "ltu"
$ws.Range("D2").Value = "ltt"
$ws.Range("C3").Value = "acu"
$ws.Range("D3").Value = "act"

# --- Denominator header, added last --------------------------------------
$ws.Range("E1").Value = "Denominator"

# --- Formatting, matching the look of the other sheets -------------------
$ws.Range("A1:C1").Font.Bold = $true
$ws.Range("E1").Font.Bold = $true
$ws.Range("C2:D3").HorizontalAlignment = -4108

$ws.Columns.Item(1).AutoFit()
$ws.Columns.Item(2).AutoFit()
$ws.Columns.Item(3).AutoFit()
$ws.Columns.Item(4).AutoFit()

# --- Selections ------------------------------------------------------------
$compartments = $wb.Worksheets.Item("Compartments")
$compartments.Range("B6:B7").Select()

$transitionParameters = $wb.Worksheets.Item("Transition Parameters")
$transitionParameters.Range("A10").Select()

# New sheet ends up active, with E1 selected (where the user was about to
# start filling in the denominator column).
$ws.Activate()
$ws.Range("E1").Select()
